# May Deskcount.xlsx update
# - Update "Include in Occupancy Calculation" (column F) from Yes to No
#   for rows 39, 43, 48, 49
# - Update Deskcount (column C) for row 45 (Melbourne) from 30 to 32
# - Update the saved selection/scroll position on the active sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

$ws.Range("F39").Value = "No"
$ws.Range("F43").Value = "No"
$ws.Range("C45").Value = 32
$ws.Range("F48").Value = "No"
$ws.Range("F49").Value = "No"

# Reflect the scrolled/selected state captured in the saved file.
$ws.Activate()
$ws.Range("C46").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 2
